$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.331.51'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.549.16'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.99'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.480'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.79'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.770.87'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.546.73'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.324.30'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.61'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.77'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.06'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0674'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.90'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.17'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.71'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.22'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.383.06'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.47'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.91'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.509'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.776'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0454'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.37'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.90'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.683.67'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.868'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.49'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.17'

$ws.Range("E2").Value = '  -0.89%  '
$ws.Range("E3").Value = '  -1.93%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("E9").Value = '  -1.96%  '
$ws.Range("E10").Value = '  -1.67%  '
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("E16").Value = '  -2.03%  '
$ws.Range("E18").Value = '  -1.90%  '
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  -2.99%  '
$ws.Range("E24").Value = '  -2.70%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -2.13%  '
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  -3.44%  '
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("E31").Value = '  -4.46%  '
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E35").Value = '  +1.30%  '
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("E39").Value = '  -2.67%  '
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("E43").Value = '  -2.16%  '
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("E45").Value = '  -1.91%  '
$ws.Range("E46").Value = '  -2.18%  '
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("E48").Value = '  -9.79%  '
$ws.Range("E49").Value = '  -1.13%  '
$ws.Range("E50").Value = '  +4.88%  '
$ws.Range("E51").Value = '  -0.29%  '
